$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently follows the
#    H1 title at the top of the document.
# ---------------------------------------------------------------------------
$metaRange = $d.Content
$metaFound = $metaRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($metaFound) {
    # Re-wrap the hit in a fresh Range and grow it to cover the whole
    # paragraph (including its trailing paragraph mark) before deleting it.
    $metaParaRange = $d.Range($metaRange.Start, $metaRange.End)
    [void]$metaParaRange.Expand(4)
    [void]$metaParaRange.Delete()
}

# ---------------------------------------------------------------------------
# 2) At the end of the document, the final paragraph currently holds the
#    (now unwanted) image-generation prompt text in italics. Split it into
#    two paragraphs:
#      - a new bold paragraph reading "Play Buffalo Mania Free: Review of
#        Game Features"
#      - the existing (italic) paragraph, now holding the meta-description
#        text that used to live near the top of the document.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$replaceRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$newXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Buffalo Mania Free: Review of Game Features</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the exciting gameplay features of Buffalo Mania, including three Jackpots and Slippery Wild Feature. Play for free and read the full review.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$replaceRange.InsertXML($newXml)
